# Reorders the data rows (2..30) of the sheet so that each target row ends up
# with the column D, H, J, K, L, M, N, O, P values that originally belonged to
# another (source) row. Columns A, B, C, E, F, G, I, Q, R are constant across
# all rows, so they do not need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (values are read from the
# ORIGINAL worksheet state before any writes happen).
$rowMap = @{
    2  = 14
    3  = 15
    4  = 9
    5  = 25
    6  = 2
    7  = 6
    8  = 30
    9  = 18
    10 = 19
    11 = 21
    12 = 27
    13 = 3
    14 = 29
    15 = 23
    16 = 7
    17 = 24
    18 = 20
    19 = 5
    20 = 26
    21 = 13
    22 = 11
    23 = 16
    24 = 12
    25 = 8
    26 = 4
    27 = 17
    28 = 10
    29 = 22
    30 = 28
}

# Columns whose values move together with the row during the reshuffle.
$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P")

# 1. Snapshot the original values for the relevant columns, for every data row.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# 2. Write the snapshot values back out according to the row mapping.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $srcVals[$c]
    }
}
